$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value2 = $aCell.Value2 + 2

    $eCell = $ws.Cells.Item($r, 5)
    $old = [string]$eCell.Value2
    if ($old -match '^(\d{2})\.(\d{2})\.(\d{4})(\d*)$') {
        $day = [int]$matches[1] + 2
        $newVal = "{0:D2}.{1}.{2}{3}" -f $day, $matches[2], $matches[3], $matches[4]
        $eCell.Value2 = $newVal
    }
}
